# Weekly CompStat sheet refresh: new report week + updated crime counts.
# (commit: "New crime data collected")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Header: bump the report Volume/Number and the covered week dates.
#    These are rich-text shared strings in the source file (several runs
#    sharing one identical font); re-assigning the whole cell text is the
#    COM-level equivalent and reproduces the exact same visible text.
# ------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "Volume 31   Number  37"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  9/9/2024  Through  9/15/2024"

# ------------------------------------------------------------------
# 2. A handful of data cells flip from a numeric 0/low count to the sheets
#    "no data" text markers ("0" or "***.*") used throughout the template.
#    Typing a numeric-looking literal ("0") gets parsed back as a number, so
#    those cells are first forced to Text format, then restored to the shared
#    look of the other text cells (copy/paste the formatting from C14, which
#    already holds the literal text "0") so the cell style matches the sheets
#    existing text-cell style rather than gaining a one-off duplicate style.
# ------------------------------------------------------------------
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "0"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(15, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 3).Value = "0"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(25, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = "0"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(27, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "0"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(28, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(28, 4).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "0"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(29, 7).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "0"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(30, 7).PasteSpecial(-4122)  # xlPasteFormats

# "***.*" is not numeric-looking, so it is stored as text as-is; only the
# cell style needs to be switched to the shared text-cell style.
$ws.Cells.Item(28, 5).Value = "***.*"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(28, 5).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(29, 8).Value = "***.*"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(29, 8).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(30, 8).Value = "***.*"
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(30, 8).PasteSpecial(-4122)  # xlPasteFormats

# ------------------------------------------------------------------
# 3. Refreshed weekly/28-day/YTD/2-year crime counts and their % changes.
# ------------------------------------------------------------------
$ws.Cells.Item(14, 13).Value = 100
$ws.Cells.Item(15, 4).Value = 2
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 7).Value = 4
$ws.Cells.Item(15, 8).Value = -75
$ws.Cells.Item(15, 10).Value = 23
$ws.Cells.Item(15, 11).Value = -52.173913043478
$ws.Cells.Item(15, 12).Value = -50
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = -50
$ws.Cells.Item(16, 6).Value = 10
$ws.Cells.Item(16, 7).Value = 18
$ws.Cells.Item(16, 8).Value = -44.444444444444
$ws.Cells.Item(16, 9).Value = 161
$ws.Cells.Item(16, 10).Value = 146
$ws.Cells.Item(16, 11).Value = 10.273972602739
$ws.Cells.Item(16, 12).Value = 4.545454545454
$ws.Cells.Item(16, 13).Value = -19.5
$ws.Cells.Item(17, 6).Value = 40
$ws.Cells.Item(17, 7).Value = 38
$ws.Cells.Item(17, 8).Value = 5.263157894736
$ws.Cells.Item(17, 9).Value = 348
$ws.Cells.Item(17, 10).Value = 337
$ws.Cells.Item(17, 11).Value = 3.264094955489
$ws.Cells.Item(17, 12).Value = -8.661417322834
$ws.Cells.Item(17, 13).Value = 45.606694560669
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 7).Value = 19
$ws.Cells.Item(18, 8).Value = -52.631578947368
$ws.Cells.Item(18, 9).Value = 89
$ws.Cells.Item(18, 10).Value = 136
$ws.Cells.Item(18, 11).Value = -34.558823529411
$ws.Cells.Item(18, 12).Value = -39.041095890411
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(19, 3).Value = 12
$ws.Cells.Item(19, 4).Value = 12
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 29
$ws.Cells.Item(19, 7).Value = 45
$ws.Cells.Item(19, 8).Value = -35.555555555555
$ws.Cells.Item(19, 9).Value = 247
$ws.Cells.Item(19, 10).Value = 319
$ws.Cells.Item(19, 11).Value = -22.570532915360
$ws.Cells.Item(19, 12).Value = -0.803212851405
$ws.Cells.Item(19, 13).Value = 17.061611374407
$ws.Cells.Item(20, 3).Value = 3
$ws.Cells.Item(20, 4).Value = 6
$ws.Cells.Item(20, 5).Value = -50
$ws.Cells.Item(20, 6).Value = 6
$ws.Cells.Item(20, 7).Value = 23
$ws.Cells.Item(20, 8).Value = -73.913043478260
$ws.Cells.Item(20, 9).Value = 63
$ws.Cells.Item(20, 10).Value = 98
$ws.Cells.Item(20, 11).Value = -35.714285714285
$ws.Cells.Item(20, 12).Value = 18.867924528301
$ws.Cells.Item(20, 13).Value = 65.789473684210
$ws.Cells.Item(21, 3).Value = 30
$ws.Cells.Item(21, 4).Value = 37
$ws.Cells.Item(21, 5).Value = -18.918918918918
$ws.Cells.Item(21, 6).Value = 95
$ws.Cells.Item(21, 7).Value = 147
$ws.Cells.Item(21, 8).Value = -35.374149659863
$ws.Cells.Item(21, 9).Value = 927
$ws.Cells.Item(21, 10).Value = 1066
$ws.Cells.Item(21, 11).Value = -13.039399624765
$ws.Cells.Item(21, 12).Value = -8.308605341246
$ws.Cells.Item(21, 13).Value = 15.442092154420
$ws.Cells.Item(23, 3).Value = 3
$ws.Cells.Item(23, 4).Value = 11
$ws.Cells.Item(23, 5).Value = -72.727272727272
$ws.Cells.Item(23, 6).Value = 17
$ws.Cells.Item(23, 7).Value = 23
$ws.Cells.Item(23, 8).Value = -26.086956521739
$ws.Cells.Item(23, 9).Value = 194
$ws.Cells.Item(23, 10).Value = 175
$ws.Cells.Item(23, 11).Value = 10.857142857142
$ws.Cells.Item(23, 12).Value = 21.25
$ws.Cells.Item(23, 13).Value = 67.241379310344
$ws.Cells.Item(24, 3).Value = 16
$ws.Cells.Item(24, 4).Value = 15
$ws.Cells.Item(24, 5).Value = 6.666666666666
$ws.Cells.Item(24, 6).Value = 59
$ws.Cells.Item(24, 7).Value = 65
$ws.Cells.Item(24, 8).Value = -9.230769230769
$ws.Cells.Item(24, 9).Value = 588
$ws.Cells.Item(24, 10).Value = 620
$ws.Cells.Item(24, 11).Value = -5.161290322580
$ws.Cells.Item(24, 12).Value = -0.675675675675
$ws.Cells.Item(24, 13).Value = 27.272727272727
$ws.Cells.Item(25, 5).Value = -100
$ws.Cells.Item(25, 6).Value = 5
$ws.Cells.Item(25, 8).Value = -64.285714285714
$ws.Cells.Item(25, 10).Value = 124
$ws.Cells.Item(25, 11).Value = -14.516129032258
$ws.Cells.Item(25, 12).Value = -31.612903225806
$ws.Cells.Item(26, 3).Value = 14
$ws.Cells.Item(26, 4).Value = 15
$ws.Cells.Item(26, 5).Value = -6.666666666666
$ws.Cells.Item(26, 6).Value = 53
$ws.Cells.Item(26, 7).Value = 45
$ws.Cells.Item(26, 8).Value = 17.777777777777
$ws.Cells.Item(26, 9).Value = 511
$ws.Cells.Item(26, 10).Value = 422
$ws.Cells.Item(26, 11).Value = 21.090047393364
$ws.Cells.Item(26, 12).Value = 18.837209302325
$ws.Cells.Item(26, 13).Value = -17.580645161290
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 7).Value = 5
$ws.Cells.Item(27, 8).Value = -60
$ws.Cells.Item(27, 10).Value = 37
$ws.Cells.Item(27, 11).Value = -51.351351351351
$ws.Cells.Item(27, 12).Value = -51.351351351351
$ws.Cells.Item(28, 6).Value = 8
$ws.Cells.Item(28, 8).Value = 166.666666666667
$ws.Cells.Item(28, 9).Value = 38
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = -30.909090909090
$ws.Cells.Item(29, 13).Value = -51.612903225806
$ws.Cells.Item(29, 14).Value = -83.146067415730
$ws.Cells.Item(30, 13).Value = -53.846153846153
$ws.Cells.Item(30, 14).Value = -85.365853658536
